$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the "Run 50" column (AZ) entirely; this shifts the "Mean" column
#    (old BA) left into AZ, and the dimension collapses from BA14 to AZ14.
$ws.Range("AZ1:AZ14").EntireColumn.Delete()

# 2) Rename the first header cell from "Gen" to "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 3) Update the "MaxFES" column values (A2:A14) - fractional progress values
#    instead of generation counts.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# 4) Update the recalculated "Mean" column (now AZ, after the column delete
#    shifted it left from BA) with the new 50-run averages.
$ws.Range("AZ2").Value = 123.868251
$ws.Range("AZ3").Value = 88.85000894
$ws.Range("AZ4").Value = 29.50491734
$ws.Range("AZ5").Value = 28.9399122
$ws.Range("AZ6").Value = 28.9399122
$ws.Range("AZ7").Value = 28.9399122
$ws.Range("AZ8").Value = 28.9399122
$ws.Range("AZ9").Value = 28.9399122
$ws.Range("AZ10").Value = 28.9399122
$ws.Range("AZ11").Value = 28.9399122
$ws.Range("AZ12").Value = 28.9399122
$ws.Range("AZ13").Value = 28.9399122
$ws.Range("AZ14").Value = 28.9399122
